$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("D11").Value = 19860602
$ws.Range("E11").Value = 17096855
$ws.Range("F11").Value = 22732975
$ws.Range("G11").Value = 24547791
$ws.Range("H11").Value = 23406044
$ws.Range("I11").Value = 29179973
$ws.Range("J11").Value = 27559098
$ws.Range("K11").Value = 36380423
$ws.Range("L11").Value = 39766680
$ws.Range("M11").Value = 33209159

$ws.Range("D12").Value = -12233707
$ws.Range("E12").Value = -10609946
$ws.Range("F12").Value = -12678171
$ws.Range("G12").Value = -15123410
$ws.Range("H12").Value = -15589570
$ws.Range("I12").Value = -20016331
$ws.Range("J12").Value = -18869377
$ws.Range("K12").Value = -24474156
$ws.Range("L12").Value = -34153538
$ws.Range("M12").Value = -27559878

$ws.Range("D13").Value = 7626895
$ws.Range("E13").Value = 6486909
$ws.Range("F13").Value = 10054804
$ws.Range("G13").Value = 9424381
$ws.Range("H13").Value = 7816474
$ws.Range("I13").Value = 9163642
$ws.Range("J13").Value = 8689721
$ws.Range("K13").Value = 11906267
$ws.Range("L13").Value = 5613142
$ws.Range("M13").Value = 5649281

$ws.Range("D14").Value = -234584
$ws.Range("E14").Value = -689678
$ws.Range("F14").Value = -5860629
$ws.Range("G14").Value = -1036663
$ws.Range("H14").Value = -115827
$ws.Range("I14").Value = -803498
$ws.Range("J14").Value = -975160
$ws.Range("K14").Value = -878829
$ws.Range("L14").Value = -1768931
$ws.Range("M14").Value = -558665

$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

$ws.Range("D16").Value = -180476
$ws.Range("E16").Value = 294513
$ws.Range("F16").Value = -2135466
$ws.Range("G16").Value = -1493178
$ws.Range("H16").Value = 1107920
$ws.Range("I16").Value = -303291
$ws.Range("J16").Value = -244958
$ws.Range("K16").Value = -800537
$ws.Range("L16").Value = -904137
$ws.Range("M16").Value = -346563

$ws.Range("D17").Value = 7211835
$ws.Range("E17").Value = 6091744
$ws.Range("F17").Value = 2058709
$ws.Range("G17").Value = 6894540
$ws.Range("H17").Value = 8808567
$ws.Range("I17").Value = 8056853
$ws.Range("J17").Value = 7469603
$ws.Range("K17").Value = 10226901
$ws.Range("L17").Value = 2940074
$ws.Range("M17").Value = 4744053

$ws.Range("D18").Value = -46556
$ws.Range("E18").Value = -68328
$ws.Range("F18").Value = -182988
$ws.Range("G18").Value = -96569
$ws.Range("H18").Value = -156963
$ws.Range("I18").Value = -133396
$ws.Range("J18").Value = -338277
$ws.Range("K18").Value = -373248
$ws.Range("L18").Value = -388435
$ws.Range("M18").Value = -373906

$ws.Range("D19").Value = -332194
$ws.Range("E19").Value = 474588
$ws.Range("F19").Value = 365354
$ws.Range("G19").Value = 209833
$ws.Range("H19").Value = 423105
$ws.Range("I19").Value = 739511
$ws.Range("J19").Value = 417092
$ws.Range("K19").Value = 537775
$ws.Range("L19").Value = 2111462
$ws.Range("M19").Value = 416668

$ws.Range("D20").Value = 6833085
$ws.Range("E20").Value = 6498004
$ws.Range("F20").Value = 2241075
$ws.Range("G20").Value = 7007804
$ws.Range("H20").Value = 9074709
$ws.Range("I20").Value = 8662968
$ws.Range("J20").Value = 7548418
$ws.Range("K20").Value = 10391428
$ws.Range("L20").Value = 4663101
$ws.Range("M20").Value = 4786815

$ws.Range("D21").Value = -594497
$ws.Range("E21").Value = -700485
$ws.Range("F21").Value = 542770
$ws.Range("G21").Value = -627665
$ws.Range("H21").Value = -974391
$ws.Range("I21").Value = -1048494
$ws.Range("J21").Value = 1835718
$ws.Range("K21").Value = -977239
$ws.Range("L21").Value = 33044
$ws.Range("M21").Value = -505149

$ws.Range("D22").Value = 6238588
$ws.Range("E22").Value = 5797519
$ws.Range("F22").Value = 2783845
$ws.Range("G22").Value = 6380139
$ws.Range("H22").Value = 8100318
$ws.Range("I22").Value = 7614474
$ws.Range("J22").Value = 9384136
$ws.Range("K22").Value = 9414189
$ws.Range("L22").Value = 4696145
$ws.Range("M22").Value = 4281666

$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

$ws.Range("D24").Value = 6238588
$ws.Range("E24").Value = 5797519
$ws.Range("F24").Value = 2783845
$ws.Range("G24").Value = 6380139
$ws.Range("H24").Value = 8100318
$ws.Range("I24").Value = 7614474
$ws.Range("J24").Value = 9384136
$ws.Range("K24").Value = 9414189
$ws.Range("L24").Value = 4696145
$ws.Range("M24").Value = 4281666

$ws.Range("D25").Value = 648
$ws.Range("E25").Value = 602
$ws.Range("F25").Value = 289
$ws.Range("G25").Value = 663
$ws.Range("H25").Value = 841
$ws.Range("I25").Value = 791
$ws.Range("J25").Value = 975
$ws.Range("K25").Value = 978
$ws.Range("L25").Value = 488
$ws.Range("M25").Value = 113

$ws.Range("D26").Value = 9629740
$ws.Range("E26").Value = 9629740
$ws.Range("F26").Value = 9629740
$ws.Range("G26").Value = 9629740
$ws.Range("H26").Value = 9629740
$ws.Range("I26").Value = 9629740
$ws.Range("J26").Value = 9629740
$ws.Range("K26").Value = 9629740
$ws.Range("L26").Value = 9629740
$ws.Range("M26").Value = 37965074

$ws.Range("D27").Value = 164
$ws.Range("E27").Value = 153
$ws.Range("F27").Value = 73
$ws.Range("G27").Value = 168
$ws.Range("H27").Value = 213
$ws.Range("I27").Value = 201
$ws.Range("J27").Value = 247
$ws.Range("K27").Value = 248
$ws.Range("L27").Value = 124
$ws.Range("M27").Value = 113
